# Generate Report for Handback
#
# Refreshes the CI-generated handback status report: the "c295f86c..." and
# "f734eb35..." entries (rows 4 and 5 on every sheet) re-ran through the
# pipeline, so their recorded priority flips from "ht" to "mt" and the
# various generation / handoff / handback timestamps advance.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: refreshed "Latest HO Xliff Generate Date" column.
$wsOverview.Range("G4").Value = "2016-08-14 16:22:51"
$wsOverview.Range("G5").Value = "2016-08-14 16:22:51"

# zh-cn sheet: Priority moved from "ht" to "mt".
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

# zh-cn sheet: refreshed Correspond Handoff / Handback datetimes.
$wsZhCn.Range("H4").Value = "2016-08-14 16:22:43"
$wsZhCn.Range("H5").Value = "2016-08-14 16:22:43"

$wsZhCn.Range("K4").Value = "2016-08-14 16:23:15"
$wsZhCn.Range("K5").Value = "2016-08-14 16:23:15"

# de-de sheet: Priority moved from "ht" to "mt".
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# de-de sheet: refreshed Correspond Handoff Datetime (mirrors the Overview
# sheet's "Latest HO Xliff Generate Date" - both were backed by the same
# shared string in the source workbook).
$wsDeDe.Range("H4").Value = "2016-08-14 16:22:51"
$wsDeDe.Range("H5").Value = "2016-08-14 16:22:51"

# de-de sheet: refreshed Correspond Handback DateTime.
$wsDeDe.Range("K4").Value = "2016-08-14 16:23:24"
$wsDeDe.Range("K5").Value = "2016-08-14 16:23:24"
